# Apply cryptos list update (GitHub Actions bot data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.638.61"
$ws.Range("E2").Value = "  -2.71%  "

$ws.Range("D3").Value = "3.743.07"
$ws.Range("E3").Value = "  -3.65%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'619.02"
$ws.Range("E5").Value = "  +2.33%  "

$ws.Range("D6").Value = "'184.86"
$ws.Range("E6").Value = "  +4.22%  "

$ws.Range("D7").Value = "3.738.92"
$ws.Range("E7").Value = "  -3.72%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").Value = "'0.729"
$ws.Range("E10").Value = "  -3.16%  "

$ws.Range("D11").Value = "'0.164"
$ws.Range("E11").Value = "  -7.93%  "

$ws.Range("D12").Value = "'58.29"
$ws.Range("E12").Value = "  +7.12%  "

$ws.Range("D13").Value = "'0.0000300"
$ws.Range("E13").Value = "  -7.49%  "

$ws.Range("D14").Value = "'10.85"
$ws.Range("E14").Value = "  -5.69%  "

$ws.Range("D15").Value = "4.336.19"
$ws.Range("E15").Value = "  -3.90%  "

$ws.Range("D16").Value = "3.738.96"
$ws.Range("E16").Value = "  -3.71%  "

$ws.Range("D17").Value = "'19.63"
$ws.Range("E17").Value = "  -6.33%  "

$ws.Range("D18").Value = "'13.05"
$ws.Range("E18").Value = "  -6.61%  "

$ws.Range("E19").Value = "  -6.38%  "

$ws.Range("D20").Value = "'0.126"
$ws.Range("E20").Value = "  -2.05%  "

$ws.Range("D21").Value = "69.417.18"
$ws.Range("E21").Value = "  -2.82%  "

$ws.Range("D22").Value = "'417.74"
$ws.Range("E22").Value = "  -5.35%  "

$ws.Range("D23").Value = "'4.77"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("E24").Value = "  -4.61%  "

$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "'11.14"
$ws.Range("E26").Value = "  -5.98%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'12.83"
$ws.Range("E27").Value = "  -7.93%  "

$ws.Range("D28").Value = "'3.90"
$ws.Range("E28").Value = "  -3.32%  "

$ws.Range("E29").Value = "  +1.61%  "

$ws.Range("E30").Value = "  -7.94%  "

$ws.Range("D31").Value = "'33.32"
$ws.Range("E31").Value = "  -5.36%  "

$ws.Range("D32").Value = "'7.49"
$ws.Range("E32").Value = "  -15.35%  "

$ws.Range("D33").Value = "'12.66"
$ws.Range("E33").Value = "  -7.11%  "

$ws.Range("E34").Value = "  -4.99%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'44.47"
$ws.Range("E35").Value = "  -8.09%  "

$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'66.21"
$ws.Range("E36").Value = "  -4.77%  "

$ws.Range("D37").Value = "'618.20"
$ws.Range("E37").Value = "  -2.64%  "

$ws.Range("D38").Value = "0.0₃0895"
$ws.Range("E38").Value = "  -11.23%  "

$ws.Range("D39").Value = "'0.413"
$ws.Range("E39").Value = "  -5.55%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("D43").Value = "'3.11"
$ws.Range("E43").Value = "  -7.06%  "

$ws.Range("D44").Value = "'0.0447"
$ws.Range("E44").Value = "  -5.44%  "

$ws.Range("E45").Value = "  -5.92%  "

$ws.Range("D46").Value = "'2.83"
$ws.Range("E46").Value = "  -10.77%  "

$ws.Range("E47").Value = "  -9.59%  "

$ws.Range("E48").Value = "  -5.33%  "

$ws.Range("D49").Value = "2.811.08"
$ws.Range("E49").Value = "  -3.44%  "

$ws.Range("E50").Value = "  -4.45%  "

$ws.Range("D51").Value = "'3.13"
$ws.Range("E51").Value = "  -3.00%  "

